# Add a "distance" column (column N) to Sheet1, matching every data row
# with a placeholder value of -1 (to be filled in later by the data
# pipeline). Mirrors the existing header-row formatting by copying it
# from an existing header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row

# Header cell
$ws.Range("N1").Value = "distance"

# Match the look of the other header cells (J1/K1 use the same header
# style already present in the sheet).
[void]$ws.Range("J1").Copy()
[void]$ws.Range("N1").PasteSpecial(-4122)

# Fill the new column with the placeholder distance value (-1) for
# every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 14).Value = -1
}

# Reflect the new column in the frozen-pane selection / scroll state.
[void]$ws.Range("N2:N77").Select()
